# Insert a new data row at row 42 (shifts existing rows 42-70 down to 43-71,
# carrying their formatting with them), then populate the new row with the
# new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(42).Insert()

$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44777
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 100112026
$ws.Range("G42").Value = "Haba"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 65
$ws.Range("K42").Value = 17000
$ws.Range("L42").Value = 17000
$ws.Range("M42").Value = 17000
$ws.Range("N42").Value = "$/saco 25 kilos"
$ws.Range("O42").Value = "Provincia de Limarí"
$ws.Range("P42").Value = 680
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
